{"js": "// Replace emoji-prefixed heading labels with bold markdown-style prefixes,\n// and bump the footer's \"Last edited\" timestamp.\n\nconst headingReplacements = [\n  { from: \"\ud83c\udfaf Nesting Complexity Levels\", to: \"**Goal:** Nesting Complexity Levels\" },\n  { from: \"\ud83d\udea8 Worst-Case Scenario Documentation\", to: \"**Note:** Worst-Case Scenario Documentation\" },\n  { from: \"\ud83c\udfaf Resolution Patterns & Rules\", to: \"**Goal:** Resolution Patterns & Rules\" },\n  { from: \"\ud83d\udd0d Edge Case Testing Matrix\", to: \"**Analysis:** Edge Case Testing Matrix\" },\n  { from: \"\ud83d\udca1 Implementation Guidelines\", to: \"**Tip:** Implementation Guidelines\" },\n  { from: \"\ud83c\udfaf Testing & Validation\", to: \"**Goal:** Testing & Validation\" },\n];\n\nconst body = context.document.body;\n\nfor (const { from, to } of headingReplacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Update the footer's \"Last edited\" date in every section.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const footer = sections.items[i].getFooter(Word.HeaderFooterType.primary);\n  const dateResults = footer.search(\"Last edited: 2025-09-07 12:53\", { matchCase: true });\n  dateResults.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < dateResults.items.length; j++) {\n    dateResults.items[j].insertText(\"Last edited: 2025-09-12 17:37\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace emoji-prefixed heading labels with bold markdown-style prefixes,\n# and bump the footer's \"Last edited\" timestamp.\n\n$d = $word.ActiveDocument\n\n$headingReplacements = @(\n    @{ From = \"\ud83c\udfaf Nesting Complexity Levels\"; To = \"**Goal:** Nesting Complexity Levels\" },\n    @{ From = \"\ud83d\udea8 Worst-Case Scenario Documentation\"; To = \"**Note:** Worst-Case Scenario Documentation\" },\n    @{ From = \"\ud83c\udfaf Resolution Patterns & Rules\"; To = \"**Goal:** Resolution Patterns & Rules\" },\n    @{ From = \"\ud83d\udd0d Edge Case Testing Matrix\"; To = \"**Analysis:** Edge Case Testing Matrix\" },\n    @{ From = \"\ud83d\udca1 Implementation Guidelines\"; To = \"**Tip:** Implementation Guidelines\" },\n    @{ From = \"\ud83c\udfaf Testing & Validation\"; To = \"**Goal:** Testing & Validation\" }\n)\n\nforeach ($pair in $headingReplacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format, ReplaceWith,\n    # Replace(2=wdReplaceAll)\n    $rng.Find.Execute($pair.From, $true, $false, $false, $false, $false, $true, 1, $false, $pair.To, 2)\n}\n\n# Update the \"Last edited\" date in every section footer.\nfor ($i = 1; $i -le $d.Sections.Count; $i++) {\n    $footer = $d.Sections($i).Footers(1)\n    $frng = $footer.Range\n    $frng.Find.ClearFormatting()\n    $frng.Find.Replacement.ClearFormatting()\n    $frng.Find.Execute(\"Last edited: 2025-09-07 12:53\", $true, $false, $false, $false, $false, $true, 1, $false, \"Last edited: 2025-09-12 17:37\", 2)\n}\n"}
